$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1538.0952

$ws.Range("H72").Value = 1538.0952

$ws.Range("H137").Value = 1903.1538
$ws.Range("I137").Value = 1569.8572
$ws.Range("K137").Value = 4709.571599999999
$ws.Range("M137").Value = -2159.571599999999

$ws.Range("H138").Value = 38465196
$ws.Range("I138").Value = 200004420
$ws.Range("J138").Value = 3476.9524
$ws.Range("K138").Value = 600013260
$ws.Range("L138").Value = 10430.8572
$ws.Range("M138").Value = -600008120
$ws.Range("N138").Value = -20710.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 6366.6665
$ws.Range("I33").Value = 6366.6665
$ws.Range("K33").Value = 6366.6665
$ws.Range("M33").Value = -6037.6665

$ws.Range("H74").Value = 58825548
$ws.Range("I74").Value = 142857820
$ws.Range("J74").Value = 2949.4
$ws.Range("K74").Value = 142857820
$ws.Range("L74").Value = 2949.4
$ws.Range("M74").Value = -142856946
$ws.Range("N74").Value = -4697.4

$ws.Range("H77").Value = 58825548
$ws.Range("I77").Value = 142857820
$ws.Range("J77").Value = 2949.4
$ws.Range("K77").Value = 714289100
$ws.Range("L77").Value = 14747
$ws.Range("M77").Value = -714284732
$ws.Range("N77").Value = -23483

$ws.Range("H122").Value = 2051.6667
$ws.Range("I122").Value = 1392.4286
$ws.Range("K122").Value = 4177.2858
$ws.Range("M122").Value = -1727.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 619
$ws.Range("I94").Value = 401.05884
$ws.Range("J94").Value = 1030.6666
$ws.Range("K94").Value = 401.05884
$ws.Range("L94").Value = 1030.6666
$ws.Range("M94").Value = 49.94116000000002
$ws.Range("N94").Value = -1932.6666

$ws.Range("H99").Value = 1701.5
$ws.Range("I99").Value = 1347.1111
$ws.Range("K99").Value = 1347.1111
$ws.Range("M99").Value = 150.8888999999999

$ws.Range("H105").Value = 2084692.1
$ws.Range("I105").Value = 1242.7333
$ws.Range("J105").Value = 5557108
$ws.Range("K105").Value = 1242.7333
$ws.Range("L105").Value = 5557108
$ws.Range("M105").Value = 504.2666999999999
$ws.Range("N105").Value = -5560602

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = $null

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").Value = $null

$ws.Range("H31").Value = 3428.6667
$ws.Range("I31").Value = 2456.8333
$ws.Range("J31").Value = 4400.5
$ws.Range("K31").Value = 2456.8333
$ws.Range("L31").Value = 4400.5
$ws.Range("M31").Value = -2161.8333
$ws.Range("N31").Value = -4990.5

$ws.Range("H34").Value = 3428.6667
$ws.Range("I34").Value = 2456.8333
$ws.Range("J34").Value = 4400.5
$ws.Range("K34").Value = 2456.8333
$ws.Range("L34").Value = 4400.5
$ws.Range("M34").Value = -2254.8333
$ws.Range("N34").Value = -4804.5

$ws.Range("H94").Value = 3061.25
$ws.Range("I94").Value = 550
$ws.Range("J94").Value = 4568
$ws.Range("K94").Value = 550
$ws.Range("L94").Value = 4568
$ws.Range("M94").Value = -99
$ws.Range("N94").Value = -5470

$ws.Range("H122").Value = 1491.4
$ws.Range("I122").Value = 1485.7778
$ws.Range("J122").Value = 1499.8334
$ws.Range("K122").Value = 4457.3334
$ws.Range("L122").Value = 4499.5002
$ws.Range("M122").Value = -2007.3334
$ws.Range("N122").Value = -9399.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 600
$ws.Range("I47").Value = 400
$ws.Range("J47").Value = 1000
$ws.Range("K47").Value = 1200
$ws.Range("L47").Value = 3000
$ws.Range("M47").Value = -769
$ws.Range("N47").Value = -3862

$ws.Range("H95").Value = 5027
$ws.Range("J95").Value = 5027
$ws.Range("L95").Value = 15081
$ws.Range("N95").Value = -19199

$ws.Range("H129").Value = 244083.9
$ws.Range("I129").Value = 630
$ws.Range("J129").Value = 511883.2
$ws.Range("K129").Value = 1890
$ws.Range("L129").Value = 1535649.6
$ws.Range("M129").Value = 3110
$ws.Range("N129").Value = -1545649.6

$ws.Range("H131").Value = 734.1799999999999
$ws.Range("J131").Value = 742.14435
$ws.Range("L131").Value = 2226.43305
$ws.Range("N131").Value = -12306.43305

$ws.Range("H134").Value = 2738.4783
$ws.Range("I134").Value = 1639.3334
$ws.Range("J134").Value = 4799.375
$ws.Range("K134").Value = 4918.0002
$ws.Range("L134").Value = 14398.125
$ws.Range("M134").Value = 151.9997999999996
$ws.Range("N134").Value = -24538.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4818364.5
$ws.Range("I70").Value = 3371.111
$ws.Range("J70").Value = 15652100
$ws.Range("K70").Value = 3371.111
$ws.Range("L70").Value = 15652100
$ws.Range("M70").Value = -3101.111
$ws.Range("N70").Value = -15652640

$ws.Range("H73").Value = 4818364.5
$ws.Range("I73").Value = 3371.111
$ws.Range("J73").Value = 15652100
$ws.Range("K73").Value = 3371.111
$ws.Range("L73").Value = 15652100
$ws.Range("M73").Value = -2435.111
$ws.Range("N73").Value = -15653972

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1010.8571
$ws.Range("I46").Value = 962.7222
$ws.Range("J46").Value = 1299.6666
$ws.Range("K46").Value = 962.7222
$ws.Range("L46").Value = 1299.6666
$ws.Range("M46").Value = -774.7222
$ws.Range("N46").Value = -1675.6666

$ws.Range("H68").Value = 1642.1428
$ws.Range("I68").Value = 1532.5
$ws.Range("J68").Value = 2300
$ws.Range("K68").Value = 1532.5
$ws.Range("L68").Value = 2300
$ws.Range("M68").Value = -783.5
$ws.Range("N68").Value = -3798

$ws.Range("H71").Value = 1642.1428
$ws.Range("I71").Value = 1532.5
$ws.Range("J71").Value = 2300
$ws.Range("K71").Value = 7662.5
$ws.Range("L71").Value = 11500
$ws.Range("M71").Value = -3918.5
$ws.Range("N71").Value = -18988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 8400
$ws.Range("J33").Value = 8400
$ws.Range("L33").Value = 8400
$ws.Range("N33").Value = -8900

$ws.Range("H36").Value = 8400
$ws.Range("J36").Value = 8400
$ws.Range("L36").Value = 8400
$ws.Range("N36").Value = -8900

$ws.Range("H107").Value = 3031190.2
$ws.Range("I107").Value = 554.875
$ws.Range("J107").Value = 6494773.5
$ws.Range("K107").Value = 1664.625
$ws.Range("L107").Value = 19484320.5
$ws.Range("M107").Value = 255.375
$ws.Range("N107").Value = -19488160.5

$ws.Range("H132").Value = 1785.68
$ws.Range("I132").Value = 1258.1538
$ws.Range("K132").Value = 3774.4614
$ws.Range("M132").Value = -1244.4614
